# Timesheet stuff and minor changes:
#  - add a "Total 8.5" label under the Hours column
#  - tidy up the one-off style on the three "Worked on Question Tool" cells
#    so they share the regular bordered cell style instead of a redundant
#    duplicate font/style
#  - move the active selection down to D12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Total 8.5" label under the Hours column (D11)
$ws.Range("D11").Value = "Total 8.5"

# Normalize the style on E7/E9/E10 back to the plain bordered style
# by re-asserting the font that already matches the regular cell style —
# this makes the engine reuse that existing style instead of the
# now-redundant duplicate font/style these cells used to carry.
$ws.Range("E7").Font.Name = "Calibri"
$ws.Range("E7").Font.Size = 11

$ws.Range("E9").Font.Name = "Calibri"
$ws.Range("E9").Font.Size = 11

$ws.Range("E10").Font.Name = "Calibri"
$ws.Range("E10").Font.Size = 11

# Move the active selection to D12
$ws.Range("D12").Select() | Out-Null
